# Generate Report for Handback
#
# This script mirrors a "handback" localization report generation:
#  - Overview / per-locale "Status" columns move from "Ready for handoff"
#    to "Handed back: in sync with en-US"
#  - The per-locale sheets (zh-cn, de-de) get their "Latest Target File",
#    "Latest Handback File" and "Latest Handback DateTime" columns filled
#    in with the handback xlf files + timestamps, and the Source File Name
#    hyperlink is duplicated into the (now populated) Target File cell.
#  - A handful of columns are widened to comfortably fit the new content.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdUrlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6b53056f7a2a1aa6d8c46b843201e0da660b920/e2e/"

$file1Name = "3e55b075-fe56-4878-a8e2-bd22e8829534.md"
$file2Name = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.md"

# ---------------------------------------------------------------------
# Overview sheet: update the per-locale status cells (E/F, rows 2 & 3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

# Widen the zh-cn / de-de status columns on the Overview sheet
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Per-locale sheets: zh-cn and de-de
# ---------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; Suffix = "zh-cn"; HandbackDate = "2016-08-28 12:48:30" },
    @{ Sheet = "de-de"; Suffix = "de-de"; HandbackDate = "2016-08-28 12:48:37" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)
    $suffix = $locale.Suffix

    $handbackFile1 = "3e55b075-fe56-4878-a8e2-bd22e8829534.7c46fd869bf0173b1a5dbb5b11cc31785398ccab.$suffix.xlf"
    $handbackFile2 = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.350b4e1ab16e1615031d3b12cf4c507ead14f7a2.$suffix.xlf"

    # Status column now reflects the handback
    $ws.Range("C2").Value = $statusNew
    $ws.Range("C3").Value = $statusNew

    # Latest Handback File / DateTime
    $ws.Range("J2").Value = $handbackFile1
    $ws.Range("K2").Value = $locale.HandbackDate

    $ws.Range("J3").Value = $handbackFile2
    $ws.Range("K3").Value = $locale.HandbackDate

    # Latest Target File gets the source file name, hyperlinked just like
    # column A's "Source File Name".
    $ws.Range("I2").Value = $file1Name
    $ws.Range("I3").Value = $file2Name

    # Rebuild the hyperlinks collection so the new "Latest Target File"
    # links (I2/I3) are interleaved with the existing "Source File Name"
    # links (A2/A3) in reading order: A2, I2, A3, I3.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "$mdUrlBase$file1Name", "", "", $file1Name) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I2"), "$mdUrlBase$file1Name", "", "", $file1Name) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), "$mdUrlBase$file2Name", "", "", $file2Name) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), "$mdUrlBase$file2Name", "", "", $file2Name) | Out-Null

    # Widen the columns that now hold longer content
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668   # Status
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664   # Latest Target File
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664  # Latest Handback File
}
